# "loading and cau template"
# Replace the "GPU tracking" accuracy table (A1:G4) with the MRC-APE
# ablation table (A1:E7): new headers, a 2-row header stack (Method/APE
# with Pre./Rec./F1 sub-columns, plus a delta-F1 column), and five data
# rows (MRC-APE (Ours) + four ablations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old table's extent (A1:G4) first.
$ws.Range("A1:G4").Clear()

# The score columns (B:E) hold numeric-looking text ("41.83", "38.20", …)
# that must stay literal text (so trailing zeros like the "0" in "38.20"
# survive) instead of being auto-coerced into floating point numbers.
# Mark the block as Text before writing, then drop the now-unneeded
# number-format style once the literal strings are in place.
$scores = $ws.Range("B3:E7")
$scores.NumberFormat = "@"

# ---- Row 1: top header ----
$ws.Range("A1").Value = "Method"
$ws.Range("B1").Value = "APE"
$ws.Range("E1").Value = "∆(F1)"
$ws.Range("B1:D1").Merge()

# ---- Row 2: sub header ----
$ws.Range("B2").Value = "Pre."
$ws.Range("C2").Value = "Rec."
$ws.Range("D2").Value = "F1"

# ---- Row 3: MRC-APE (Ours) ----
$ws.Range("A3").Value = "MRC-APE (Ours)"
$ws.Range("B3").Value = "41.83"
$ws.Range("C3").Value = "38.17"
$ws.Range("D3").Value = "39.92"
$ws.Range("E3").Value = "-"

# ---- Row 4: w/o Db -> Da ----
$ws.Range("A4").Value = "w/o Db → Da"
$ws.Range("B4").Value = "49.47"
$ws.Range("C4").Value = "31.33"
$ws.Range("D4").Value = "38.36"
$ws.Range("E4").Value = "1.56"

# ---- Row 5: w/o Da -> Db ----
$ws.Range("A5").Value = "w/o Da → Db"
$ws.Range("B5").Value = "46.68"
$ws.Range("C5").Value = "26.02"
$ws.Range("D5").Value = "33.41"
$ws.Range("E5").Value = "6.51"

# ---- Row 6: w/o LSTM ----
$ws.Range("A6").Value = "w/o LSTM"
$ws.Range("B6").Value = "44.98"
$ws.Range("C6").Value = "34.51"
$ws.Range("D6").Value = "39.06"
$ws.Range("E6").Value = "0.86"

# ---- Row 7: w/o GA ----
$ws.Range("A7").Value = "w/o GA"
$ws.Range("B7").Value = "38.20"
$ws.Range("C7").Value = "30.66"
$ws.Range("D7").Value = "34.02"
$ws.Range("E7").Value = "5.90"

# Values are written; the Text number-format was only needed to stop the
# numeric-looking strings above from being parsed as floats, so drop it
# again now to leave the cells on the sheet's default style.
$scores.ClearFormats()

# Every populated cell in this table is its own single-cell merge region,
# matching the workbook's existing per-cell table convention (B1:D1 above
# is the one genuine multi-cell merge).
$ws.Range("A1").Merge()
$ws.Range("E1").Merge()
$ws.Range("B2").Merge()
$ws.Range("C2").Merge()
$ws.Range("D2").Merge()
$ws.Range("A3").Merge()
$ws.Range("B3").Merge()
$ws.Range("C3").Merge()
$ws.Range("D3").Merge()
$ws.Range("E3").Merge()
$ws.Range("A4").Merge()
$ws.Range("B4").Merge()
$ws.Range("C4").Merge()
$ws.Range("D4").Merge()
$ws.Range("E4").Merge()
$ws.Range("A5").Merge()
$ws.Range("B5").Merge()
$ws.Range("C5").Merge()
$ws.Range("D5").Merge()
$ws.Range("E5").Merge()
$ws.Range("A6").Merge()
$ws.Range("B6").Merge()
$ws.Range("C6").Merge()
$ws.Range("D6").Merge()
$ws.Range("E6").Merge()
$ws.Range("A7").Merge()
$ws.Range("B7").Merge()
$ws.Range("C7").Merge()
$ws.Range("D7").Merge()
$ws.Range("E7").Merge()
